# Updated cryptos list: price (D) and 1h volume-change (E) columns refreshed; rows
# 44-46 and 51 also changed coin identity (name/link swap with updated price data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value even when it looks numeric,
    # matching the original inlineStr cell type/style (no format residue).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "92.501.63"
Set-TextValue $ws.Range("E2") "  -2.04%  "
Set-TextValue $ws.Range("D3") "3.396.17"
Set-TextValue $ws.Range("E3") "  -0.79%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "229.51"
Set-TextValue $ws.Range("E5") "  -3.27%  "
Set-TextValue $ws.Range("D6") "616.34"
Set-TextValue $ws.Range("E6") "  -4.07%  "
Set-TextValue $ws.Range("E7") "  -5.30%  "
Set-TextValue $ws.Range("E8") "  -4.19%  "
Set-TextValue $ws.Range("E9") "  +0.06%  "
Set-TextValue $ws.Range("D10") "0.955"
Set-TextValue $ws.Range("E10") "  -1.67%  "
Set-TextValue $ws.Range("D11") "3.394.95"
Set-TextValue $ws.Range("E11") "  -0.75%  "
Set-TextValue $ws.Range("D12") "42.65"
Set-TextValue $ws.Range("E12") "  +1.98%  "
Set-TextValue $ws.Range("E13") "  -1.51%  "
Set-TextValue $ws.Range("D14") "6.21"
Set-TextValue $ws.Range("E14") "  +0.00%  "
Set-TextValue $ws.Range("D15") "4.044.24"
Set-TextValue $ws.Range("E15") "  -0.63%  "
Set-TextValue $ws.Range("D16") "92.400.43"
Set-TextValue $ws.Range("E16") "  -1.88%  "
Set-TextValue $ws.Range("E17") "  -2.91%  "
Set-TextValue $ws.Range("D18") "8.02"
Set-TextValue $ws.Range("E18") "  -3.76%  "
Set-TextValue $ws.Range("D19") "3.398.13"
Set-TextValue $ws.Range("E19") "  -0.73%  "
Set-TextValue $ws.Range("D20") "17.73"
Set-TextValue $ws.Range("E20") "  +1.35%  "
Set-TextValue $ws.Range("D21") "11.41"
Set-TextValue $ws.Range("E21") "  -1.94%  "
Set-TextValue $ws.Range("D22") "493.55"
Set-TextValue $ws.Range("E22") "  -1.17%  "
Set-TextValue $ws.Range("E23") "  +2.05%  "
Set-TextValue $ws.Range("D24") "0.432"
Set-TextValue $ws.Range("E24") "  -12.82%  "
Set-TextValue $ws.Range("D25") "6.49"
Set-TextValue $ws.Range("E25") "  -0.29%  "
Set-TextValue $ws.Range("D26") "0.0000182"
Set-TextValue $ws.Range("E26") "  -5.42%  "
Set-TextValue $ws.Range("D27") "90.09"
Set-TextValue $ws.Range("E27") "  -4.16%  "
Set-TextValue $ws.Range("D28") "11.88"
Set-TextValue $ws.Range("E28") "  -0.55%  "
Set-TextValue $ws.Range("D29") "3.577.55"
Set-TextValue $ws.Range("E29") "  -0.80%  "
Set-TextValue $ws.Range("E30") "  +0.00%  "
Set-TextValue $ws.Range("D31") "11.17"
Set-TextValue $ws.Range("E31") "  -4.95%  "
Set-TextValue $ws.Range("E32") "  -2.43%  "
Set-TextValue $ws.Range("D33") "0.133"
Set-TextValue $ws.Range("E33") "  -3.31%  "
Set-TextValue $ws.Range("D34") "0.987"
Set-TextValue $ws.Range("E34") "  -1.14%  "
Set-TextValue $ws.Range("D35") "0.171"
Set-TextValue $ws.Range("E35") "  -4.50%  "
Set-TextValue $ws.Range("D36") "29.38"
Set-TextValue $ws.Range("E36") "  -0.69%  "
Set-TextValue $ws.Range("D37") "0.537"
Set-TextValue $ws.Range("E37") "  -2.33%  "
Set-TextValue $ws.Range("D38") "551.10"
Set-TextValue $ws.Range("E38") "  +0.03%  "
Set-TextValue $ws.Range("D39") "7.41"
Set-TextValue $ws.Range("E39") "  -3.06%  "
Set-TextValue $ws.Range("E40") "  -0.03%  "
Set-TextValue $ws.Range("E41") "  -1.13%  "
Set-TextValue $ws.Range("E42") "  -5.13%  "
Set-TextValue $ws.Range("D43") "0.906"
Set-TextValue $ws.Range("E43") "  +0.60%  "
Set-TextValue $ws.Range("D47") "5.42"
Set-TextValue $ws.Range("E47") "  -3.22%  "
Set-TextValue $ws.Range("D48") "0.0402"
Set-TextValue $ws.Range("E48") "  -2.06%  "
Set-TextValue $ws.Range("D49") "53.20"
Set-TextValue $ws.Range("E49") "  -4.26%  "
Set-TextValue $ws.Range("E50") "  -4.49%  "

# Rows 44-46: coin ordering/identities shifted (MantraDAO <-> ImmutableX <-> WhiteBITCoin)
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D44") "1.71"
Set-TextValue $ws.Range("E44") "  -0.86%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D45") "23.61"
Set-TextValue $ws.Range("E45") "  -1.98%  "

$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D46") "3.67"
Set-TextValue $ws.Range("E46") "  +1.46%  "

# Row 51: Fantom replaced by Cosmos
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D51") "7.90"
Set-TextValue $ws.Range("E51") "  -1.33%  "
